$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting existing rows 103:139 down to 104:140
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new price-record data
$ws.Cells.Item(103, 1).Value = 4
$ws.Cells.Item(103, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(103, 3).Value = "Los Lagos"
$ws.Cells.Item(103, 4).Value = 45027
$ws.Cells.Item(103, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(103, 5).Value = 10
$ws.Cells.Item(103, 6).Value = 100112031
$ws.Cells.Item(103, 7).Value = "Poroto verde"
$ws.Cells.Item(103, 8).Value = "Magnum"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 40
$ws.Cells.Item(103, 11).Value = 25000
$ws.Cells.Item(103, 12).Value = 25000
$ws.Cells.Item(103, 13).Value = 25000
$ws.Cells.Item(103, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(103, 15).Value = "Región Metropolitana"
$ws.Cells.Item(103, 16).Value = 1000
$ws.Cells.Item(103, 17).Value = 25
$ws.Cells.Item(103, 18).Value = "Hortaliza"
